$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "#aggiornamento 15/09" note cell (J9). ClearContents also
# drops the now-unused shared string and shrinks the row's span/dimension.
$ws.Range("J9").ClearContents() | Out-Null

# Give column J the width that used to live on column M (closest value the
# host lets us set - engine quantizes ColumnWidth to 1/6-character steps).
$ws.Columns.Item(10).ColumnWidth = 9.585

# Refresh the weekly progress figures (rows 66:131): new report date plus
# updated ticket counts / average handle times.
$ws.Range("B66").Value2 = 45918
$ws.Range("C66").Value2 = 79
$ws.Range("D66").Value2 = 33.0434898785425
$ws.Range("B67").Value2 = 45918
$ws.Range("C67").Value2 = 44
$ws.Range("D67").Value2 = 43.1817792207792
$ws.Range("B68").Value2 = 45918
$ws.Range("C68").Value2 = 96
$ws.Range("D68").Value2 = 27.4786069303082
$ws.Range("B69").Value2 = 45918
$ws.Range("C69").Value2 = 56
$ws.Range("D69").Value2 = 39.7565939039409
$ws.Range("B70").Value2 = 45918
$ws.Range("C70").Value2 = 104
$ws.Range("D70").Value2 = 29.1326935483871
$ws.Range("B71").Value2 = 45918
$ws.Range("C71").Value2 = 64
$ws.Range("D71").Value2 = 43.5968314864582
$ws.Range("B72").Value2 = 45918
$ws.Range("C72").Value2 = 120
$ws.Range("D72").Value2 = 32.0599444444445
$ws.Range("B73").Value2 = 45918
$ws.Range("C73").Value2 = 72
$ws.Range("D73").Value2 = 33.7079425287356
$ws.Range("B74").Value2 = 45918
$ws.Range("C74").Value2 = 96
$ws.Range("D74").Value2 = 50.7345554037006
$ws.Range("B75").Value2 = 45918
$ws.Range("C75").Value2 = 104
$ws.Range("D75").Value2 = 55.4522929336223
$ws.Range("B76").Value2 = 45918
$ws.Range("C76").Value2 = 52
$ws.Range("D76").Value2 = 33.8694145299145
$ws.Range("B77").Value2 = 45918
$ws.Range("C77").Value2 = 64
$ws.Range("D77").Value2 = 48.8747588243867
$ws.Range("B78").Value2 = 45918
$ws.Range("C78").Value2 = 96
$ws.Range("D78").Value2 = 29.1326935483871
$ws.Range("B79").Value2 = 45918
$ws.Range("C79").Value2 = 96
$ws.Range("D79").Value2 = 30.5464905626781
$ws.Range("B80").Value2 = 45918
$ws.Range("C80").Value2 = 104
$ws.Range("D80").Value2 = 31.6699326923077
$ws.Range("B81").Value2 = 45918
$ws.Range("C81").Value2 = 8
$ws.Range("D81").Value2 = 48.4532523364486
$ws.Range("B82").Value2 = 45918
$ws.Range("C82").Value2 = 80
$ws.Range("D82").Value2 = 33.0434898785425
$ws.Range("B83").Value2 = 45918
$ws.Range("C83").Value2 = 104
$ws.Range("D83").Value2 = 29.7360729166667
$ws.Range("B84").Value2 = 45918
$ws.Range("C84").Value2 = 104
$ws.Range("D84").Value2 = 48.1536690961653
$ws.Range("B85").Value2 = 45918
$ws.Range("C85").Value2 = 88
$ws.Range("D85").Value2 = 33.4594291428269
$ws.Range("B86").Value2 = 45918
$ws.Range("C86").Value2 = 104
$ws.Range("D86").Value2 = 52.9118425676887
$ws.Range("B87").Value2 = 45918
$ws.Range("C87").Value2 = 72
$ws.Range("D87").Value2 = 52.1261519795658
$ws.Range("B88").Value2 = 45918
$ws.Range("C88").Value2 = 104
$ws.Range("D88").Value2 = 39.7934126565649
$ws.Range("B89").Value2 = 45918
$ws.Range("C89").Value2 = 100
$ws.Range("D89").Value2 = 44.5242949464397
$ws.Range("B90").Value2 = 45918
$ws.Range("C90").Value2 = 88
$ws.Range("D90").Value2 = 112.065143292683
$ws.Range("B91").Value2 = 45918
$ws.Range("C91").Value2 = 88
$ws.Range("D91").Value2 = 32.0225772727273
$ws.Range("B92").Value2 = 45918
$ws.Range("C92").Value2 = 86
$ws.Range("D92").Value2 = 25.0018244274809
$ws.Range("B93").Value2 = 45918
$ws.Range("C93").Value2 = 108
$ws.Range("D93").Value2 = 44.8153288271811
$ws.Range("B94").Value2 = 45918
$ws.Range("C94").Value2 = 96
$ws.Range("D94").Value2 = 40.3412567327236
$ws.Range("B95").Value2 = 45918
$ws.Range("C95").Value2 = 112
$ws.Range("D95").Value2 = 30.8888623511905
$ws.Range("B96").Value2 = 45918
$ws.Range("C96").Value2 = 80
$ws.Range("D96").Value2 = 58.1200511108775
$ws.Range("B97").Value2 = 45918
$ws.Range("C97").Value2 = 104
$ws.Range("D97").Value2 = 65.3275448717949
$ws.Range("B98").Value2 = 45918
$ws.Range("C98").Value2 = 40
$ws.Range("D98").Value2 = 26.656775
$ws.Range("B99").Value2 = 45918
$ws.Range("C99").Value2 = 104
$ws.Range("D99").Value2 = 55.8999651698649
$ws.Range("B100").Value2 = 45918
$ws.Range("C100").Value2 = 88
$ws.Range("D100").Value2 = 33.8209034090909
$ws.Range("B101").Value2 = 45918
$ws.Range("C101").Value2 = 88
$ws.Range("D101").Value2 = 112.065143292683
$ws.Range("B102").Value2 = 45918
$ws.Range("C102").Value2 = 96
$ws.Range("D102").Value2 = 25.7496293363019
$ws.Range("B103").Value2 = 45918
$ws.Range("C103").Value2 = 92
$ws.Range("D103").Value2 = 33.9754833333333
$ws.Range("B104").Value2 = 45918
$ws.Range("C104").Value2 = 80
$ws.Range("D104").Value2 = 34.1102432692308
$ws.Range("B105").Value2 = 45918
$ws.Range("C105").Value2 = 120
$ws.Range("D105").Value2 = 33.789925
$ws.Range("B106").Value2 = 45918
$ws.Range("C106").Value2 = 72
$ws.Range("D106").Value2 = 77.6300070974675
$ws.Range("B107").Value2 = 45918
$ws.Range("C107").Value2 = 68
$ws.Range("D107").Value2 = 33.0434898785425
$ws.Range("B108").Value2 = 45918
$ws.Range("C108").Value2 = 80
$ws.Range("D108").Value2 = 32.7152125
$ws.Range("B109").Value2 = 45918
$ws.Range("C109").Value2 = 104
$ws.Range("D109").Value2 = 65.3275448717949
$ws.Range("B110").Value2 = 45918
$ws.Range("C110").Value2 = 80
$ws.Range("D110").Value2 = 65.3275448717949
$ws.Range("B111").Value2 = 45918
$ws.Range("C111").Value2 = 120
$ws.Range("D111").Value2 = 31.6427083333333
$ws.Range("B112").Value2 = 45918
$ws.Range("C112").Value2 = 80
$ws.Range("D112").Value2 = 27.0203472222222
$ws.Range("B113").Value2 = 45918
$ws.Range("C113").Value2 = 120
$ws.Range("D113").Value2 = 31.7142086111111
$ws.Range("B114").Value2 = 45918
$ws.Range("C114").Value2 = 120
$ws.Range("D114").Value2 = 28.9230021043771
$ws.Range("B115").Value2 = 45918
$ws.Range("C115").Value2 = 72
$ws.Range("D115").Value2 = 30.5797861467237
$ws.Range("B116").Value2 = 45918
$ws.Range("C116").Value2 = 104
$ws.Range("D116").Value2 = 35.1611832314881
$ws.Range("B117").Value2 = 45918
$ws.Range("C117").Value2 = 80
$ws.Range("D117").Value2 = 33.1502161675824
$ws.Range("B118").Value2 = 45918
$ws.Range("C118").Value2 = 120
$ws.Range("D118").Value2 = 28.4317083333333
$ws.Range("B119").Value2 = 45918
$ws.Range("C119").Value2 = 40
$ws.Range("B120").Value2 = 45918
$ws.Range("C120").Value2 = 104
$ws.Range("D120").Value2 = 29.1326935483871
$ws.Range("B121").Value2 = 45918
$ws.Range("C121").Value2 = 88
$ws.Range("D121").Value2 = 66.4255282991819
$ws.Range("B122").Value2 = 45918
$ws.Range("C122").Value2 = 104
$ws.Range("D122").Value2 = 44.736744413734
$ws.Range("B123").Value2 = 45918
$ws.Range("C123").Value2 = 96
$ws.Range("D123").Value2 = 45.7111822916667
$ws.Range("B124").Value2 = 45918
$ws.Range("C124").Value2 = 112
$ws.Range("D124").Value2 = 31.1100080891331
$ws.Range("B125").Value2 = 45918
$ws.Range("C125").Value2 = 48
$ws.Range("D125").Value2 = 38.8561849856322
$ws.Range("B126").Value2 = 45918
$ws.Range("C126").Value2 = 72
$ws.Range("D126").Value2 = 69.7983911124451
$ws.Range("B127").Value2 = 45918
$ws.Range("C127").Value2 = 120
$ws.Range("D127").Value2 = 33.789925
$ws.Range("B128").Value2 = 45918
$ws.Range("C128").Value2 = 104
$ws.Range("D128").Value2 = 142.196326923077
$ws.Range("B129").Value2 = 45918
$ws.Range("C129").Value2 = 72
$ws.Range("D129").Value2 = 33.9754833333333
$ws.Range("B130").Value2 = 45918
$ws.Range("C130").Value2 = 32
$ws.Range("D130").Value2 = 32.7152125
$ws.Range("B131").Value2 = 45918
$ws.Range("C131").Value2 = 32
$ws.Range("D131").Value2 = 28.685

# Restore the author's last on-screen scroll position / selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 94
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I113").Select() | Out-Null
